{"js": "const body = context.document.body;\nbody.insertParagraph(\"hello world\", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$d.Content.InsertParagraphAfter()\n"}
